# Migrate Excel creation to XSSF backend (#254)
# -----------------------------------------------------------------
# Reproduces the user-visible / content-level edits that accompanied
# the template resave:
#   1. Typo fix in the shared string used by cell A48.
#   2. Selection moved to row 10 (whole row selected, active cell A10).
#   3. Header/footer font style name corrected ("Normal" -> "Regular").
# -----------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix the typo "gedankeloses" -> "gedankenloses"
$ws.Range("A48").Value = "Das Denken der Gedanken ist ein gedankenloses Denken"

# 2. Update the active selection to row 10
$ws.Range("A10").EntireRow.Select()

# 3. Header / footer: "Times New Roman,Normal" -> "Times New Roman,Regular"
$ps = $ws.PageSetup
$ps.CenterHeader = '&"Times New Roman,Regular"&12&A'
$ps.CenterFooter = '&"Times New Roman,Regular"&12Página &P'
